$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "Tue Nov 12 18:05:28 EST 2024"
$ws.Range("B3").Value = "Tue Nov 12 18:05:42 EST 2024"
$ws.Range("B4").Value = "Tue Nov 12 18:05:55 EST 2024"
